# Rename the worksheet: "Sheet2" -> "Sheet1"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet1"

# Update the regression statistics in F28:N28 (Future time perspective row).
# Downstream formulas in row 31 (TEXT(...) summaries) recalc automatically.
$ws.Range("F28").Value = 0.0266646890059696
$ws.Range("G28").Value = 0.481614740523214
$ws.Range("H28").Value = 0.684798605676998
$ws.Range("I28").Value = 0.0364789545084105
$ws.Range("J28").Value = 0.65908334361017
$ws.Range("K28").Value = 0.744846241533684
$ws.Range("L28").Value = 0.0362927227264766
$ws.Range("M28").Value = 0.655714152587062
$ws.Range("N28").Value = 0.743764728371043
